$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Low)
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5764.04
$ws.Range("D2").Value = 2933.7
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 70

# Row 3 (mid-Low)
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 560.96
$ws.Range("D3").Value = 2584.57
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 9002.940000000001

# Row 4 (Middle)
$ws.Range("B4").Value = 1212
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 94.72
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = 5369.71

# Row 5 (mid-High)
$ws.Range("B5").Value = 9647.67
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 1807
$ws.Range("F5").Value = 282.68

# Row 6 (High)
$ws.Range("B6").Value = 4434.33
$ws.Range("E6").Value = 1142
$ws.Range("F6").ClearContents()
